# "removes NAs from recapture table"
#
# Row 2 of the Recaptures-EDIQuery sheet was a dummy/example data row where
# every text column held the placeholder string "NA" (and the numeric
# column M2 held 0). The edit clears that placeholder data out of columns
# A-P entirely, and also blanks out the (still border/fill-styled) Q:T
# cells, leaving their formatting in place but with no value. Once nothing
# in the sheet references the shared string "NA" anymore, it naturally
# drops out of the shared-strings table on save.
#
# The sheet's view is also reset: the former scroll/selection position
# (topLeftCell P1, selection T2) is replaced by a selection on B7 with no
# special scroll offset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the placeholder values in A2:P2 entirely (no leftover cell records).
$ws.Range("A2:P2").ClearContents() | Out-Null

# Q2:T2 keep their existing style (border/fill) but lose their "NA" value.
$ws.Range("Q2:T2").ClearContents() | Out-Null

# Move the selection/active cell to B7 (this also clears the old
# top-left-cell scroll anchor that pinned the view at column P).
$ws.Range("B7").Select() | Out-Null
